# Add missing "Admin document upload validation checks" (DOC-03) test case
# to the Documents section of the tracker, right after DOC-02.
# This shifts ARC-01..UI-01 down by one row (rows 15-24 -> 16-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 15, pushing existing rows 15-24 down to 16-25.
$ws.Rows("15").Insert()

# Populate the new row with the DOC-03 test case.
$ws.Range("A15").Value = "DOC-03"
$ws.Range("B15").Value = "Documents"
$ws.Range("C15").Value = "Admin document upload validation checks"
$ws.Range("D15").Value = "Admin"
$ws.Range("E15").Value = "Try unsupported extension or oversize file on admin upload form"
$ws.Range("F15").Value = "Specific validation error appears and file is not added"
$ws.Range("G15").Value = "NOT RUN"
